# Update the "Solvers" column (B) and "Count" column (C) for four partner
# rows to reflect the newly tracked solver names, per the commit:
# "added ability to track all partners for solvers".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Firefly Innovations -> now matched with AHSA Platform
$ws.Range("B29").Value = "None,AHSA Platform"
$ws.Range("C29").Value = 1

# Kevin Przybocki -> now matched with AHSA Platform and Jute-based biodegradable PPE
$ws.Range("B38").Value = "None,AHSA Platform,Jute-based biodegradable PPE"
$ws.Range("C38").Value = 2

# The Kamath Family Foundation -> now matched with AHSA Platform
$ws.Range("B69").Value = "None,AHSA Platform"
$ws.Range("C69").Value = 1

# Usizo Advisory Solutions -> now matched with AHSA Platform
$ws.Range("B76").Value = "None,AHSA Platform"
$ws.Range("C76").Value = 1
